$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The blog article window shown on row 7 advances by one:
#   I7: ser 81 -> ser 82
#   E7: ser 82 -> ser 83
#   C7: ser 83 -> ser 84 (new article 84 is now live)
# D7 (the meetup card) is left untouched.

$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 82"
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 83"
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 84"
